# thêm cột "Mã CQT cấp" cho xuất khẩu excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D. Excel shifts the old D..Z data
# (and the merged A1/A2 title banners) right to E..AA, and the new
# column inherits column C's cell formatting.
$ws.Columns("D:D").Insert()

# The insert does not carry over C's explicit width, so copy it across.
$ws.Columns("D:D").ColumnWidth = $ws.Columns("C:C").ColumnWidth

# Header text for the newly inserted column.
$ws.Cells.Item(4, 4).Value = "Mã CQT cấp"

# Matches the active-cell selection recorded after the edit.
$ws.Range("D5").Select()
